# Add the new row (row 4) with the branch note in column A, as introduced
# by this commit ("add row by alaaeddin"). This appends a new shared
# string and grows the sheet's used range / dimension to A1:B4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "This row was added by branch Alaaeddin"

# Column A needs to widen (AutoFit) so the new, longer text fits — mirrors
# the bestFit column-width recalculation Excel performs for the column.
$ws.Columns.Item(1).ColumnWidth = 36.666666666666664
